# Traditional model results: add F1-score columns and rename LogR -> LogitR
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has columns:
#   B=Model, C=Sensitivity (w/o resampling), D=Specificity (w/o resampling),
#   E=Sensitivity (w/ resampling), F=Specificity (w/ resampling)
#
# Target layout:
#   B=Model, C=Sensitivity, D=Specificity, E=F1-score, F=(blank spacer),
#   G=Sensitivity (w/ resampling), H=Specificity (w/ resampling), I=F1-Score (w/ resampling)

# Step 1: shift the existing "w/ resampling" columns (E,F) two slots to the
# right (to G,H) by inserting two new blank columns at E,F. Insert() copies
# formatting from the columns being pushed, so header formatting is preserved.
$ws.Columns("E:F").Insert()

# Step 2: rename the "w/o resampling" headers (now plain names)
$ws.Range("C1").Value = "Sensitivity"
$ws.Range("D1").Value = "Specificity"

# Step 3: add the new F1-score headers
$ws.Range("E1").Value = "F1-score"

# I1 is a brand new column outside the original range, so it does not
# inherit the bold/bordered header formatting the way the shifted columns
# did. Copy the header formatting from H1 (an existing header cell) first,
# then set its text.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1").PasteSpecial(-4122) | Out-Null
$ws.Range("I1").Value = "F1-Score (w/ resampling)"

# Step 4: fill in the new F1-score values (without resampling) in column E
$ws.Range("E2").Value = 0.8288832234665614
$ws.Range("E3").Value = 0.8350927033463117
$ws.Range("E4").Value = 0.7979078887464386

# Step 5: fill in the new F1-score values (with resampling) in column I
$ws.Range("I2").Value = 0.8418941199624707
$ws.Range("I3").Value = 0.9749928291636436
$ws.Range("I4").Value = 0.7375805834715173

# Step 6: rename model "LogR" -> "LogitR"
$ws.Range("B4").Value = "LogitR"
